# daily auto push: 2026-01-14 02:32 UTC
#
# The source data table (A:D = 日付/曜日/時刻/ランキング) gets two more
# "2026/01/14" sample rows appended to that date's block, which live
# before the (unrelated, already-present) "2026/12/29" block that follows
# it in the sheet. Concretely: insert two blank rows at 633:634 (pushing
# the existing row 633 and everything below it down by two), then fill
# those two new rows with the new readings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 633; this shifts the old
# rows 633..674 down to 635..676 and grows the used range accordingly.
$ws.Rows("633:634").Insert()

# Row 633: 2026/01/14, 水, 8, 201
# The leading "'" forces text so Excel doesn't reinterpret the
# "yyyy/mm/dd" string as a real date (matching the other date cells in
# column A, which are stored as plain text).
$ws.Range("A633").Value = "'2026/01/14"
$ws.Range("B633").Value = "水"
$ws.Range("C633").Value = 8
$ws.Range("D633").Value = 201

# Row 634: 2026/01/14, 水, 9, 201
$ws.Range("A634").Value = "'2026/01/14"
$ws.Range("B634").Value = "水"
$ws.Range("C634").Value = 9
$ws.Range("D634").Value = 201

# Drop the quote-prefix/text formatting picked up from the apostrophe
# trick above so the new cells end up with the same (default/no-style)
# formatting as every other row in the table.
$ws.Range("A633:D634").ClearFormats()
